# "radar plot early draft"
#
# The 2020-election rows (12 = Trump/Republican, 13 = Biden/Democrat) had
# their pop_vote (F), pop_vote_pc (G) and electoral_vote (H) figures swapped
# between the two rows. Put the correct figures back with each candidate:
#   Row 12 (Trump)  -> pop_vote 74216154, pop_vote_pc 46.9, electoral_vote 232
#   Row 13 (Biden)  -> pop_vote 81268924, pop_vote_pc 51.3, electoral_vote 306

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F12").Value = 74216154
$ws.Range("G12").Value = 46.9
$ws.Range("H12").Value = 232

$ws.Range("F13").Value = 81268924
$ws.Range("G13").Value = 51.3
$ws.Range("H13").Value = 306

# Leave the cursor where the author last clicked while eyeballing the
# swapped figures.
$ws.Range("G14").Select()
